$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" date column (C2:C6) from 2023-11-13 (45243) to 2023-11-14 (45244)
for ($row = 2; $row -le 6; $row++) {
    $ws.Cells.Item($row, 3).Value = 45244
}
